{"js": "// The document's last paragraph is a \"ListParagraph\" list item that reads\n// \"...simultaneously.\" and is immediately followed (inside the same\n// paragraph) by the \"_GoBack\" bookmark.\n//\n// The target edit splits that single paragraph into three list paragraphs,\n// leaving the bookmark anchored to the new, trailing (now empty) paragraph:\n//   1) \"...simultaneously.\"                                            (existing, trimmed)\n//   2) \"Finding outline asset in order to distinguish each cell of the table\" (new)\n//   3) \"\" - empty paragraph that now holds the _GoBack bookmark          (new)\n//\n// Locate the bookmark and insert text containing paragraph-mark characters\n// (\"\\r\") immediately before it. This mirrors what happens in Word when a\n// user places the cursor right before the bookmark and presses Enter twice\n// while typing the new line: the newly created paragraphs inherit the\n// paragraph style (ListParagraph) and numbering (ilvl/numId) of the\n// paragraph being split, and the bookmark itself is left untouched, ending\n// up alone in the final paragraph.\nconst bookmarkRange = context.document.getBookmarkRange(\"_GoBack\");\n\nbookmarkRange.insertText(\n  \"\\rFinding outline asset in order to distinguish each cell of the table\\r\",\n  \"Before\"\n);\n\nawait context.sync();\n", "ps1": "# The document's last paragraph is a \"ListParagraph\" list item that reads\n# \"...simultaneously.\" and is immediately followed (inside the same\n# paragraph) by the \"_GoBack\" bookmark.\n#\n# The target edit splits that single paragraph into three list paragraphs,\n# leaving the bookmark anchored to the new, trailing (now empty) paragraph:\n#   1) \"...simultaneously.\"                                                  (existing, trimmed)\n#   2) \"Finding outline asset in order to distinguish each cell of the table\" (new)\n#   3) \"\"  - empty paragraph that now holds the _GoBack bookmark             (new)\n\n$d = $word.ActiveDocument\n\n# Locate the (hidden) \"_GoBack\" bookmark and collapse its range to a caret\n# right at its position (still \"inside\" the original paragraph, immediately\n# after the trailing \".\" run).\n$bm = $d.Bookmarks.Item('_GoBack')\n$r = $bm.Range\n$r.Collapse(1)  # wdCollapseStart\n\n# Insert the new line of text preceded and followed by a paragraph mark. This\n# mirrors a user placing the cursor right before the bookmark and pressing\n# Enter, typing the new line, then pressing Enter again: the two freshly\n# created paragraphs inherit the paragraph style (ListParagraph) and\n# numbering (ilvl/numId) of the paragraph being split, and the bookmark is\n# left untouched, ending up alone in the final paragraph.\n$r.InsertBefore(\"`rFinding outline asset in order to distinguish each cell of the table`r\")\n"}
